$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Formatted Number to Text:" / "Blank Text:" row (old row 40, C40) holds a
# shared-string cell with no value, which reads back as blank. Make sure it
# stays genuinely blank before everything shifts down.
$ws.Range("C40").ClearContents()

# Insert a new row above the old "@ format to Number:" row (old row 34) for
# the new "Percentage Text to Number:" example. This pushes every row from
# 34 downward down by one, matching rows 35-43 in the target layout.
$ws.Rows("34").Insert()

# New row: percentage text gets parsed to a number and displayed with a
# percentage number format (e.g. "55.12%").
$ws.Range("B34").Value = "Percentage Text to Number:"
$ws.Range("C34").Value = 0.5512
$ws.Range("C34").NumberFormat = "0.00%"

# Column B needs to be a bit wider to fit the new, longer label.
$ws.Columns("B").ColumnWidth = 25.15
